$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: split the "...counts from 1 to 10." run so "10." is wrapped
# in a gramStart/gramEnd proofErr pair (paragraph "Problem:  What finger
# will the girl land on if she counts from 1 to 10.").
# ---------------------------------------------------------------------
$target1 = "What finger will the girl land on if she counts from 1 to 10."
$p1 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.Contains($target1)) {
        $p1 = $cand
        break
    }
}

$xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' `
    + '<w:r w:rsidRPr="00267887"><w:t>Problem</w:t></w:r>' `
    + '<w:r><w:t>:</w:t></w:r>' `
    + '<w:r w:rsidR="00C36E54"><w:t xml:space="preserve">  What finger will the girl land on if she counts from 1 to </w:t></w:r>' `
    + '<w:proofErr w:type="gramStart"/>' `
    + '<w:r w:rsidR="00C36E54"><w:t>10.</w:t></w:r>' `
    + '<w:proofErr w:type="gramEnd"/>' `
    + '</w:p>'

$p1.Range.InsertXML($xml1)

# ---------------------------------------------------------------------
# Change 2: the "C.   If she counts from 1 to 100 = stop on ring finger
# again." paragraph becomes an auto-numbered ListParagraph item (losing
# its manual "C.   " prefix and its ind/left formatting), and is
# followed by four blank paragraphs and a new explanatory paragraph
# ("Why? ..."), which now carries the _GoBack bookmark that used to sit
# on the "ring finger again." paragraph.
# ---------------------------------------------------------------------
$target2 = "ring finger again."
$p2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.Contains($target2)) {
        $p2 = $cand
        break
    }
}

$whyText = "Why? Counting from to 10 lands on first finger.  In order to figure out what she would land on from 1 to 100, you have to multiply 10 x 10 = 100.  Count ten fingers more from where you left on 10 and you will land where you should be for 100 instead of actually counting it all out."

$xml2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00D04725" w:rsidRPr="00267887" w:rsidRDefault="00652B7D" w:rsidP="00D04725">' `
    + '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' `
    + '<w:r><w:t xml:space="preserve">If she counts from 1 to 100 = stop on </w:t></w:r>' `
    + '<w:r w:rsidR="00985813"><w:t>ring finger again.</w:t></w:r>' `
    + '</w:p>' `
    + '<w:p/><w:p/><w:p/><w:p/>' `
    + '<w:p><w:r><w:t>' + $whyText + '</w:t></w:r>' `
    + '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$p2.Range.InsertXML($xml2)
